$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the Price column (D) as Text first so values like "1.642.83" or "1.003"
# are stored as literal strings (matching the source data) rather than being
# auto-coerced into numbers by Excel's input parser.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.961.18'
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").Value = '1.640.02'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '214.94'
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("D6").Value = '0.5067'
$ws.Range("E6").Value = '  +0.75%  '
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '0.2560'
$ws.Range("E8").Value = '  -0.48%  '
$ws.Range("D9").Value = '0.06366'
$ws.Range("E9").Value = '  -0.20%  '
$ws.Range("D10").Value = '19.48'
$ws.Range("E10").Value = '  -0.45%  '
$ws.Range("D11").Value = '0.07764'
$ws.Range("E11").Value = '  +0.30%  '
$ws.Range("D12").Value = '4.286'
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").Value = '1.651.25'
$ws.Range("E13").Value = '  +0.68%  '
$ws.Range("D14").Value = '0.5448'
$ws.Range("E14").Value = '  -0.11%  '
$ws.Range("D15").Value = '0.0₅7821'
$ws.Range("E15").Value = '  -0.92%  '
$ws.Range("D16").Value = '64.33'
$ws.Range("E16").Value = '  +0.40%  '
$ws.Range("D17").Value = '26.010.59'
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("D18").Value = '1.002'
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("D19").Value = '197.79'
$ws.Range("E19").Value = '  -2.31%  '
$ws.Range("D20").Value = '4.434'
$ws.Range("E20").Value = '  +0.87%  '
$ws.Range("D21").Value = '9.954'
$ws.Range("E21").Value = '  +0.82%  '
$ws.Range("D22").Value = '6.045'
$ws.Range("E22").Value = '  +1.22%  '
$ws.Range("D23").Value = '1.004'
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("D24").Value = '1.894'
$ws.Range("E24").Value = '  +1.22%  '
$ws.Range("D25").Value = '141.26'
$ws.Range("E25").Value = '  +0.41%  '
$ws.Range("E26").Value = '  +3.59%  '
$ws.Range("D27").Value = '6.874'
$ws.Range("E27").Value = '  +1.19%  '
$ws.Range("D28").Value = '15.72'
$ws.Range("E28").Value = '  +0.18%  '
$ws.Range("D29").Value = '1.236'
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("D30").Value = '0.04979'
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("D31").Value = '3.258'
$ws.Range("E31").Value = '  -0.36%  '
$ws.Range("D32").Value = '3.183'
$ws.Range("E32").Value = '  -0.21%  '
$ws.Range("D33").Value = '1.540'
$ws.Range("E33").Value = '  -0.14%  '
$ws.Range("D34").Value = '2.356'
$ws.Range("E34").Value = '  -0.29%  '
$ws.Range("D35").Value = '0.8938'
$ws.Range("E35").Value = '  +0.36%  '
$ws.Range("D36").Value = '2.578'
$ws.Range("E36").Value = '  -1.86%  '
$ws.Range("D37").Value = '1.130.22'
$ws.Range("E37").Value = '  -1.54%  '
$ws.Range("D38").Value = '0.5445'
$ws.Range("E38").Value = '  -3.43%  '
$ws.Range("E39").Value = '  -0.58%  '
$ws.Range("D40").Value = '2.546'
$ws.Range("E40").Value = '  -0.47%  '
$ws.Range("D41").Value = '1.002'
$ws.Range("E41").Value = '  +0.13%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '5.595'
$ws.Range("E42").Value = '  -1.39%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '0.8170'
$ws.Range("E43").Value = '  +1.49%  '
$ws.Range("B44").Value = 'BabyDogeCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D44").Value = '0.0₈126'
$ws.Range("E44").Value = '  +8.10%  '
$ws.Range("D45").Value = '99.66'
$ws.Range("E45").Value = '  -0.16%  '
$ws.Range("D46").Value = '1.778.47'
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("D47").Value = '0.4537'
$ws.Range("E47").Value = '  +0.05%  '
$ws.Range("E48").Value = '  -0.10%  '
$ws.Range("D49").Value = '54.77'
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").Value = '0.05071'
$ws.Range("E50").Value = '  +0.41%  '
$ws.Range("D51").Value = '1.003'
$ws.Range("E51").Value = '  +0.37%  '

# Restore the original (default) style on the Price column now that the text
# values are committed, so no stray formatting is left behind.
$ws.Range("D2:D51").Style = "Normal"

